# Updates price (D) and 1h volume change (E) figures for the cryptos
# list, and swaps the Kaspa/Celestia rows (32/33) so Celestia is listed
# first. A leading apostrophe is used on a few Price values so Excel
# keeps them as plain text instead of auto-converting them to numbers
# (which would silently drop significant trailing zeros, e.g. 102.30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.909.21'
$ws.Range("E2").Value = '  +3.80%  '

$ws.Range("D3").Value = '2.420.25'
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''315.35'
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("D6").Value = '''102.30'
$ws.Range("E6").Value = '  +5.30%  '

$ws.Range("E7").Value = '  +1.37%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = '''0.523'
$ws.Range("E9").Value = '  +7.46%  '

$ws.Range("D10").Value = '''35.31'
$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("E11").Value = '  +0.44%  '

$ws.Range("E12").Value = '  -2.38%  '

$ws.Range("D13").Value = '''18.18'
$ws.Range("E13").Value = '  -2.39%  '

$ws.Range("D14").Value = '''6.95'
$ws.Range("E14").Value = '  +1.26%  '

$ws.Range("D15").Value = '2.804.07'
$ws.Range("E15").Value = '  +1.37%  '

$ws.Range("D16").Value = '2.431.26'
$ws.Range("E16").Value = '  +1.95%  '

$ws.Range("D17").Value = '''0.834'
$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("D18").Value = '44.821.12'
$ws.Range("E18").Value = '  +3.50%  '

$ws.Range("D19").Value = '''12.19'
$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("E20").Value = '  -1.13%  '

$ws.Range("D21").Value = '0.0₃0918'
$ws.Range("E21").Value = '  +1.96%  '

$ws.Range("D22").Value = '''68.60'
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("E23").Value = '  +1.98%  '

$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("D25").Value = '''2.47'
$ws.Range("E25").Value = '  +0.37%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").Value = '''25.10'
$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("E28").Value = '  -6.65%  '

$ws.Range("D29").Value = '''9.50'
$ws.Range("E29").Value = '  +0.91%  '

$ws.Range("D30").Value = '''48.96'
$ws.Range("E30").Value = '  +1.69%  '

$ws.Range("D31").Value = '''32.62'
$ws.Range("E31").Value = '  +1.00%  '

$ws.Range("B32").Value = 'Celestia'
$ws.Range("C32").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D32").Value = '''19.72'
$ws.Range("E32").Value = '  +6.83%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '''0.123'
$ws.Range("E33").Value = '  +7.32%  '

$ws.Range("D34").Value = '''5.20'
$ws.Range("E34").Value = '  +1.03%  '

$ws.Range("E35").Value = '  +0.19%  '

$ws.Range("D36").Value = '''0.0756'
$ws.Range("E36").Value = '  +1.39%  '

$ws.Range("D37").Value = '''1.86'
$ws.Range("E37").Value = '  -0.68%  '

$ws.Range("D38").Value = '''4.41'
$ws.Range("E38").Value = '  +1.04%  '

$ws.Range("D39").Value = '''2.84'
$ws.Range("E39").Value = '  -6.78%  '

$ws.Range("D40").Value = '''121.83'
$ws.Range("E40").Value = '  -7.24%  '

$ws.Range("E41").Value = '  -3.06%  '

$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("D43").Value = '''20.70'
$ws.Range("E43").Value = '  -2.33%  '

$ws.Range("D44").Value = '''0.0289'
$ws.Range("E44").Value = '  +2.29%  '

$ws.Range("D45").Value = '1.931.85'
$ws.Range("E45").Value = '  -0.75%  '

$ws.Range("E46").Value = '  -2.45%  '

$ws.Range("E47").Value = '  +3.48%  '

$ws.Range("D48").Value = '''9.21'
$ws.Range("E48").Value = '  -1.33%  '

$ws.Range("E49").Value = '  +14.50%  '

$ws.Range("D50").Value = '''75.85'
$ws.Range("E50").Value = '  +5.30%  '

$ws.Range("D51").Value = '''53.17'
$ws.Range("E51").Value = '  +1.38%  '
